# The workbook's "Sheet" contains stock-report rows where, for a number of
# products, the two data rows belonging to that product had been recorded in
# the wrong order (row N holding what should be row N+1's figures, and vice
# versa). This swaps the data fields (everything except the running index in
# column A) back between each such pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-row data (A is just the running serial number
# and is left untouched; H:M are blank placeholder columns).
$cols = @("B", "C", "D", "E", "F", "G")

# Each tuple is (row holding the "wrong"/swapped data, its paired row).
$pairs = @(
    @(149, 150),
    @(264, 265),
    @(313, 314),
    @(316, 317),
    @(350, 352),
    @(355, 356),
    @(372, 373),
    @(382, 383),
    @(400, 401),
    @(419, 420),
    @(431, 432),
    @(457, 458),
    @(536, 537),
    @(586, 587),
    @(590, 591),
    @(593, 594),
    @(601, 602),
    @(604, 605),
    @(709, 710),
    @(715, 716),
    @(720, 721)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$row1")
        $cell2 = $ws.Range("$col$row2")

        $value1 = $cell1.Value()
        $value2 = $cell2.Value()

        $cell1.Value = $value2
        $cell2.Value = $value1
    }
}
